$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.020.91'
$ws.Range("E2").Value = '  -1.00%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.905.30'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7561'
$ws.Range("E5").Value = '  +1.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.53'
$ws.Range("E6").Value = '  -0.86%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.0000'
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("E8").Value = '  -1.91%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.55'
$ws.Range("E9").Value = '  -6.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06908'
$ws.Range("E10").Value = '  -1.12%  '
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7560'
$ws.Range("E12").Value = '  -2.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.898.28'
$ws.Range("E13").Value = '  -1.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.256'
$ws.Range("E14").Value = '  -0.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.85'
$ws.Range("E15").Value = '  -0.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.195'
$ws.Range("E16").Value = '  +5.63%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.026.37'
$ws.Range("E17").Value = '  -0.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.05'
$ws.Range("E18").Value = '  -1.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007758'
$ws.Range("E19").Value = '  -1.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.90'
$ws.Range("E20").Value = '  -3.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9999'
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.154.12'
$ws.Range("E22").Value = '  -1.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.054'
$ws.Range("E24").Value = '  +6.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.319'
$ws.Range("E25").Value = '  -1.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.79'
$ws.Range("E26").Value = '  +0.77%  '
$ws.Range("E27").Value = '  -0.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1288'
$ws.Range("E28").Value = '  +0.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.073'
$ws.Range("E29").Value = '  -3.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.347'
$ws.Range("E30").Value = '  -1.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.531'
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("E32").Value = '  -1.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.049'
$ws.Range("E33").Value = '  -0.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05397'
$ws.Range("E34").Value = '  +3.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.289'
$ws.Range("E35").Value = '  -1.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7377'
$ws.Range("E36").Value = '  -1.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.720'
$ws.Range("E37").Value = '  -1.89%  '
$ws.Range("E38").Value = '  +0.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.766'
$ws.Range("E39").Value = '  -0.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.255'
$ws.Range("E40").Value = '  -2.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4455'
$ws.Range("E41").Value = '  -0.44%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.78'
$ws.Range("E42").Value = '  -4.19%  '
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8315'
$ws.Range("E45").Value = '  -1.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.713'
$ws.Range("E46").Value = '  +0.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.64'
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.858'
$ws.Range("E48").Value = '  +0.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.053.52'
$ws.Range("E49").Value = '  -0.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.57'
$ws.Range("E50").Value = '  -1.78%  '
$ws.Range("E51").Value = '  -0.10%  '